# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master (drives every slide)
#   ppt/theme/theme2.xml  -> bound to the notes master
#
# The edit swaps the two themes' content: the slide master's theme
# ("Integral") becomes the default "Office Theme" colors, and vice versa
# for the notes master. We apply this by rewriting the 12 theme colors
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) on the reachable
# Master.Theme.ThemeColorScheme to the "Office Theme" palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Index -> RGB (stored as a COM long, 0xBBGGRR) for the standard "Office Theme"
# color scheme, in clrScheme order: dk1, lt1, dk2, lt2, accent1..6, hlink, folHlink.
$colorScheme.Item(1).RGB  = 0          # dk1      000000
$colorScheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477    # folHlink 954F72
